$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, exactly matching the refreshed
# crypto-price snapshot pulled by the scheduled GitHub Actions run.
$updates = [ordered]@{
    "D2" = "43.901.50"
    "E2" = "  +1.18%  "
    "D3" = "2.236.12"
    "E3" = "  -0.04%  "
    "E4" = "  +0.26%  "
    "D5" = "315.57"
    "E5" = "  -1.22%  "
    "D6" = "99.23"
    "E6" = "  -1.50%  "
    "D7" = "0.569"
    "E7" = "  -2.71%  "
    "E8" = "  +0.13%  "
    "D9" = "0.535"
    "E9" = "  -5.28%  "
    "D10" = "36.34"
    "E10" = "  -2.76%  "
    "D11" = "0.0821"
    "E11" = "  -2.30%  "
    "D12" = "7.36"
    "E12" = "  -4.96%  "
    "E13" = "  -2.72%  "
    "D14" = "2.578.49"
    "E14" = "  +0.12%  "
    "B15" = "WrappedEther"
    "C15" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D15" = "2.241.73"
    "E15" = "  -0.35%  "
    "B16" = "Polygon"
    "C16" = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
    "D16" = "0.839"
    "E16" = "  -3.18%  "
    "D17" = "14.00"
    "E17" = "  -2.31%  "
    "D18" = "43.801.77"
    "E18" = "  +0.98%  "
    "D19" = "12.79"
    "E19" = "  -10.01%  "
    "D20" = "0.0₃0964"
    "E20" = "  -2.87%  "
    "D21" = "6.34"
    "E21" = "  -4.47%  "
    "D22" = "64.84"
    "E22" = "  -1.23%  "
    "E23" = "  -3.64%  "
    "D24" = "233.42"
    "E24" = "  -1.49%  "
    "E25" = "  -7.12%  "
    "E26" = "  +0.38%  "
    "E27" = "  +1.42%  "
    "E28" = "  -1.22%  "
    "D29" = "36.80"
    "E29" = "  +0.80%  "
    "D30" = "6.02"
    "E30" = "  -5.77%  "
    "D31" = "157.91"
    "E31" = "  -1.59%  "
    "D32" = "19.92"
    "E32" = "  -1.98%  "
    "D33" = "0.0831"
    "E33" = "  -4.97%  "
    "E34" = "  -1.44%  "
    "D35" = "3.18"
    "E35" = "  -2.10%  "
    "E36" = "  +5.03%  "
    "E37" = "  -0.39%  "
    "E38" = "  -2.86%  "
    "D39" = "15.99"
    "E39" = "  +6.96%  "
    "D40" = "3.62"
    "E40" = "  -3.37%  "
    "D41" = "4.05"
    "E41" = "  -8.15%  "
    "E42" = "  -3.90%  "
    "E43" = "  +0.16%  "
    "D44" = "1.731.40"
    "E44" = "  -4.53%  "
    "D45" = "0.194"
    "E45" = "  -5.06%  "
    "D46" = "80.45"
    "E46" = "  -3.92%  "
    "D47" = "73.32"
    "E47" = "  -2.09%  "
    "E48" = "  -3.80%  "
    "B49" = "Stacks"
    "C49" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D49" = "1.65"
    "E49" = "  +0.40%  "
    "B50" = "Aave"
    "C50" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D50" = "101.34"
    "E50" = "  -1.55%  "
    "D51" = "56.67"
    "E51" = "  -3.95%  "
}

# Cells whose new text looks like a plain number (e.g. "315.57") must be
# forced to Text format first, otherwise Excel auto-converts the literal
# string into a numeric value instead of keeping it as text.
$textForce = @(
    "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D16",
    "D17", "D19", "D21", "D22", "D24", "D29", "D30", "D31",
    "D32", "D33", "D35", "D39", "D40", "D41", "D45", "D46",
    "D47", "D49", "D50", "D51"
)

foreach ($ref in $textForce) {
    $ws.Range($ref).NumberFormat = "@"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in $textForce) {
    $ws.Range($ref).ClearFormats()
}
